$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (day 43897) blog slots shift up: the article that was showing in C6
# ("ser: 65") moves to G6, and the newly published article 66 takes over C6.
$ws.Range("G6").Value = "type: blog`nwidth: 2`nheight: 1`nser: 65"
$ws.Range("C6").Value = "type: blog`nwidth: 2`nheight: 1`nser: 66"

# Reflect the author's last-saved selection on the sheet.
[void]$ws.Range("G6").Select()
